$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows, per the repull/push of data
$ws.Range("F2").Value = -1
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -12
$ws.Range("F7").Value = -7
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -5
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -5
$ws.Range("F15").Value = -3
